$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Logic Level FET"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "https://www.digikey.com/en/products/detail/infineon-technologies/BSS806NH6327XTSA1/2783472"

$ws.Range("C7").Value = "https://www.digikey.com/en/products/detail/bourns-inc/2211-H-RC/775363"

$ws.Range("C30").Select()
